# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" table: new snapshot timestamp, updated
# case counts for a number of countries, and re-ranked rows for the
# three countries whose totals crossed a neighbour's (Eslovaquia /
# Tailandia / Mozambique and Lituania / Eslovenia).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Timestamp banner in row 1
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 10:54"

# 2. Updated rows: country name (re-ranked where needed) + the 7 data columns
#    (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos
#    criticos, Muertes hoy, Muertes)
$rows = @(
    @{ Row=6;   Country='India';            B=3110761; C=5576; D=2338899; E=714147; F=0; G=23; H=57715 }
    @{ Row=25;  Country='Filipinas';        B=194252;  C=4686; D=132042;  E=59200;  F=0; G=13; H=3010 }
    @{ Row=26;  Country='Indonesia';        B=155412;  C=1877; D=111060;  E=37593;  F=0; G=79; H=6759 }
    @{ Row=33;  Country='Israel';           B=103274;  C=611;  D=80521;   E=21914;  F=0; G=5;  H=839 }
    @{ Row=47;  Country='Polonia';          B=62310;   C=548;  D=42448;   E=17902;  F=0; G=5;  H=1960 }
    @{ Row=63;  Country='Afganistan';       B=38054;   C=55;   D=28360;   E=8305;   F=0; G=2;  H=1389 }
    @{ Row=71;  Country='Austria';          B=25495;   C=242;  D=21657;   E=3105;   F=0; G=1;  H=733 }
    @{ Row=92;  Country='Malasia';          B=9274;    C=7;    D=8965;    E=184;    F=0; G=0;  H=125 }
    @{ Row=111; Country='Hong Kong';        B=4692;    C=9;    D=4052;    E=563;    F=0; G=0;  H=77 }
    @{ Row=120; Country='Eslovaquia';       B=3424;    C=68;   D=2153;    E=1238;   F=0; G=0;  H=33 }
    @{ Row=121; Country='Tailandia';        B=3397;    C=2;    D=3222;    E=117;    F=0; G=0;  H=58 }
    @{ Row=122; Country='Mozambique';       B=3356;    C=0;    D=1503;    E=1872;   F=0; G=0;  H=20 }
    @{ Row=130; Country='Lituania';         B=2673;    C=38;   D=1766;    E=822;    F=0; G=1;  H=85 }
    @{ Row=131; Country='Eslovenia';        B=2665;    C=14;   D=2122;    E=410;    F=0; G=2;  H=133 }
    @{ Row=133; Country='Estonia';          B=2275;    C=3;    D=2025;    E=186;    F=0; G=1;  H=64 }
    @{ Row=202; Country='Santa Lucia';      B=26;      C=0;    D=25;      E=1;      F=0; G=0;  H=0 }
    @{ Row=203; Country='Timor Oriental';   B=26;      C=0;    D=25;      E=1;      F=0; G=0;  H=0 }
    @{ Row=205; Country='Nueva Caledonia';  B=23;      C=0;    D=23;      E=0;      F=0; G=0;  H=0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Country
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}
